# Insert a new data row at row 18 (pushing existing rows 18:54 down to 19:55)
# and populate it with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44623
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100114002
$ws.Range("G18").Value = "Camote"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 18000
$ws.Range("N18").Value = "$/malla 20 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 900
$ws.Range("Q18").Value = 20
$ws.Range("R18").Value = "Hortaliza"
